$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: the blank paragraph right after
# "b)I can not say that all solutions will work in all cases..."
# becomes the "5) Choose a solution..." heading, and a brand new
# paragraph is inserted after it with the explanation text.
# ------------------------------------------------------------------
$headingIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^b\)I can not say that all solutions will work in all cases") {
        # The first blank paragraph right after "b)I can not say..." is left
        # untouched; the *second* blank paragraph is the one that becomes
        # the "5) Choose a solution..." heading.
        $headingIdx = $i + 2
        break
    }
}

$apostrophe = [char]8217
$heading = "5) Choose a solution and develop a plan to implement it."
$explanation = "For my specific case I will assume that the there is water in the river, that he is in a very rural area that does not have any bridges near him and that he can" + $apostrophe + "t swim. For him to efficiently attain his goal he needs to first tie the bag of seeds and throw it across the river, aiming at a soft turf where plant life meets the river then place the cat in the boat. He should swiftly travel across the river in the boat to the other side where he will see the bird trying to open he bag of seeds. The bag should not be badly damaged, the bird should be there since it would be attracted to the seeds and the cat should be safely getting out of the boat with the man."

$headingPara = $d.Paragraphs.Item($headingIdx)
$headingPara.Range.Text = $heading + "`r" + $explanation

# ------------------------------------------------------------------
# Edit 2: trim the run of 8 blank paragraphs after
# "(b)At least one matching pair of each color." down to 3 blank
# paragraphs (remove 5 of them).
# ------------------------------------------------------------------
$afterIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^\(b\)At least one matching pair of each color\.") {
        $afterIdx = $i + 1
        break
    }
}

for ($n = 0; $n -lt 5; $n++) {
    $d.Paragraphs.Item($afterIdx).Range.Delete()
}
